# Bernard Updated Test Results Sheet
#
# On the "Single Algos" sheet, duplicate row 31 (Bernard's "#44 - Top 5%"
# manual XGB entry) into a new row 32 for the follow-up "#50 - Top 5%"
# snapshot that used "boost deviance" instead of "boost exponential", and
# clear out the now-redundant "Wining Snapshot" label in B35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single Algos")

# Clone row 31 (values + formatting) into row 32.
$ws.Range("A31:M31").Copy($ws.Range("A32"))

# Update the Local AUC / Kaggle AUC results for the new snapshot.
$ws.Range("J32").Value = 0.86614999999999998
$ws.Range("K32").Value = 0.86812

# New Position / Notes text for the new snapshot.
$ws.Range("L32").Value = "#50 - Top 5%"
$ws.Range("M32").Value = "XGB Rnd 7641,3431,1270,8939,9101 /  boost deviance / adaboost"

# The "Wining Snapshot" label is no longer needed here.
$ws.Range("B35").ClearContents()

# Move the cursor to reflect where editing left off.
$null = $ws.Range("B33").Select()
